$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency price/volume snapshot. A new "Frax" entry is
# inserted at row 40 (shifting the following rows down by one), and the
# trailing "EOS" row is dropped, per the GitHub Actions data pull.
# Price/volume text is written with a leading apostrophe so Excel keeps
# it as literal text (e.g. '1.020', '20.00') instead of coercing it to a
# number and dropping significant trailing/format zeros.
$data = @(
    ,@('Bitcoin','https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc',"'27.916.05","'  +0.36%  ")
    ,@('Ethereum','https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth',"'1.884.21","'  +0.18%  ")
    ,@('TetherUSD','https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt',"'1.018","'  +1.60%  ")
    ,@('BNB','https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb',"'335.39","'  +0.49%  ")
    ,@('USDC','https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc',"'1.017","'  +1.48%  ")
    ,@('XRP','https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp',"'0.4685","'  -1.00%  ")
    ,@('Cardano','https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada',"'0.3907","'  -1.80%  ")
    ,@('OKB','https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb',"'46.82","'  -3.42%  ")
    ,@('Dogecoin','https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge',"'0.07954","'  -1.32%  ")
    ,@('Polygon','https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic',"'1.013","'  -1.37%  ")
    ,@('Solana','https://coinranking.com/coin/zNZHO_Sjf+solana-sol',"'21.72","'  -1.16%  ")
    ,@('WrappedEther','https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth',"'1.880.13","'  -0.74%  ")
    ,@('Polkadot','https://coinranking.com/coin/25W7FG7om+polkadot-dot',"'5.956","'  -0.29%  ")
    ,@('Chainlink','https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link',"'7.120","'  -0.97%  ")
    ,@('BinanceUSD','https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd',"'1.020","'  +1.46%  ")
    ,@('TRON','https://coinranking.com/coin/qUhEFk1I61atv+tron-trx',"'0.06787","'  +2.45%  ")
    ,@('Litecoin','https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc',"'87.43","'  +0.19%  ")
    ,@('ShibaInu','https://coinranking.com/coin/xz24e0BjL+shibainu-shib',"'0.00001047","'  -0.56%  ")
    ,@('Avalanche','https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax',"'17.02","'  -1.47%  ")
    ,@('Dai','https://coinranking.com/coin/MoTuySvg7+dai-dai',"'1.017","'  +1.54%  ")
    ,@('WrappedBTC','https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc',"'27.915.34","'  +0.19%  ")
    ,@('Uniswap','https://coinranking.com/coin/_H5FVG9iW+uniswap-uni',"'5.468","'  -0.74%  ")
    ,@('Cosmos','https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom',"'10.93","'  -1.10%  ")
    ,@('Toncoin','https://coinranking.com/coin/67YlI0K1b+toncoin-ton',"'2.363","'  +2.78%  ")
    ,@('WrappedliquidstakedEther2.0','https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth',"'2.109.59","'  -0.30%  ")
    ,@('Monero','https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr',"'159.86","'  +1.73%  ")
    ,@('EthereumClassic','https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc',"'20.00","'  -1.36%  ")
    ,@('LidoDAOToken','https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo',"'2.082","'  -1.20%  ")
    ,@('InternetComputer(DFINITY)','https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp',"'5.469","'  -2.58%  ")
    ,@('BitcoinCash','https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch',"'121.07","'  -1.47%  ")
    ,@('Stellar','https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm',"'0.09554","'  -0.18%  ")
    ,@('ImmutableX','https://coinranking.com/coin/Z96jIvLU7+immutablex-imx',"'0.9561","'  -2.23%  ")
    ,@('HuobiToken','https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht',"'3.655","'  +0.53%  ")
    ,@('Filecoin','https://coinranking.com/coin/ymQub4fuB+filecoin-fil',"'5.334","'  +0.22%  ")
    ,@('ARBITRUM','https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb',"'1.349","'  -8.09%  ")
    ,@('Hedera','https://coinranking.com/coin/jad286TjB+hedera-hbar',"'0.06110","'  +0.01%  ")
    ,@('VeChain','https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet',"'0.02245","'  -0.88%  ")
    ,@('TrustWalletToken','https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt',"'1.208","'  -1.70%  ")
    ,@('Frax','https://coinranking.com/coin/KfWtaeV1W+frax-frax',"'1.017","'  +1.52%  ")
    ,@('FraxShare','https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs',"'8.159","'  -0.75%  ")
    ,@('TheSandbox','https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand',"'0.5917","'  -1.98%  ")
    ,@('Algorand','https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo',"'0.1893","'  -1.06%  ")
    ,@('Aptos','https://coinranking.com/coin/HGYj5JCv5+aptos-apt',"'10.23","'  -1.05%  ")
    ,@('WEMIXTOKEN','https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix',"'1.270","'  +1.85%  ")
    ,@('Decentraland','https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana',"'0.5654","'  -1.20%  ")
    ,@('EnergySwap','https://coinranking.com/coin/SbWqqTui-+energyswap-ens',"'12.19","'  -1.33%  ")
    ,@('PancakeSwap','https://coinranking.com/coin/ncYFcP709+pancakeswap-cake',"'3.395","'  -0.76%  ")
    ,@('NEARProtocol','https://coinranking.com/coin/DCrsaMv68+nearprotocol-near',"'1.927","'  -0.91%  ")
    ,@('Cronos','https://coinranking.com/coin/65PHZTpmE55b+cronos-cro',"'0.06857","'  +0.47%  ")
    ,@('Quant','https://coinranking.com/coin/bauj_21eYVwso+quant-qnt',"'113.91","'  +0.12%  ")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 2).Value = $row[$j]
    }
}

Write-Output "Updated cryptos list"
